$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 is being re-purposed from a set of placeholder "a" / blank entries
# into a real shipping/order record. Several of the new values look
# numeric/date-like ("02/14/2024", "5.0", "50.0", "79935") but must be
# stored as literal text (matching the rest of the sheet, which was built
# with every cell as text) rather than being auto-coerced into a date or
# number by Excel - which would also silently swap in a new number-format
# style for the cell.
#
# Trick: write the literal text as a formula ("=""value"""), which forces a
# text/string result without touching the cell's number format, then
# collapse the formula down to its static value with a values-only paste.
# That keeps the original style index untouched while still ending up with
# a plain text cell.
function Set-LiteralText {
    param($addr, $text)

    $escaped = $text -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

Set-LiteralText "A12" "02/14/2024"
Set-LiteralText "B12" "SO240214001"
Set-LiteralText "C12" "AB"
Set-LiteralText "D12" "9999999999"
Set-LiteralText "E12" "ARTIST"
Set-LiteralText "F12" "TITLE"
Set-LiteralText "G12" "5.0"
Set-LiteralText "H12" "50.0"
Set-LiteralText "J12" "LP"
Set-LiteralText "K12" "Ashley"
Set-LiteralText "L12" "YES"
Set-LiteralText "M12" "test"
Set-LiteralText "N12" "test"
Set-LiteralText "O12" "te"
Set-LiteralText "P12" "79935"

$excel.CutCopyMode = $false
